$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 530.
# This shifts the existing rows 530-601 down to become rows 532-603,
# preserving all of their values and formatting.
$ws.Rows.Item(530).Insert()
$ws.Rows.Item(530).Insert()

# --- Populate the new row 530 ---
$ws.Cells.Item(530, 1).Value = 10
$ws.Cells.Item(530, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(530, 3).Value = "La Araucanía"
$ws.Cells.Item(530, 4).Value = 45124
$ws.Cells.Item(530, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(530, 5).Value = 9
$ws.Cells.Item(530, 6).Value = 100112040
$ws.Cells.Item(530, 7).Value = "Cilantro"
$ws.Cells.Item(530, 8).Value = "Sin especificar"
$ws.Cells.Item(530, 9).Value = "Primera"
$ws.Cells.Item(530, 10).Value = 60
$ws.Cells.Item(530, 11).Value = 4000
$ws.Cells.Item(530, 12).Value = 4000
$ws.Cells.Item(530, 13).Value = 4000
$ws.Cells.Item(530, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(530, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(530, 16).Value = 2000
$ws.Cells.Item(530, 17).Value = 2
$ws.Cells.Item(530, 18).Value = "Hortaliza"

# --- Populate the new row 531 ---
$ws.Cells.Item(531, 1).Value = 10
$ws.Cells.Item(531, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(531, 3).Value = "La Araucanía"
$ws.Cells.Item(531, 4).Value = 45124
$ws.Cells.Item(531, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(531, 5).Value = 9
$ws.Cells.Item(531, 6).Value = 100112040
$ws.Cells.Item(531, 7).Value = "Cilantro"
$ws.Cells.Item(531, 8).Value = "Sin especificar"
$ws.Cells.Item(531, 9).Value = "Primera"
$ws.Cells.Item(531, 10).Value = 100
$ws.Cells.Item(531, 11).Value = 4600
$ws.Cells.Item(531, 12).Value = 4600
$ws.Cells.Item(531, 13).Value = 4600
$ws.Cells.Item(531, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(531, 15).Value = "Región Metropolitana"
$ws.Cells.Item(531, 16).Value = 2300
$ws.Cells.Item(531, 17).Value = 2
$ws.Cells.Item(531, 18).Value = "Hortaliza"
